$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: task #100008 - Create Map Component
$ws.Range("A10").Value = "#100008"
$ws.Range("B10").Value = "Create Map Component"

# Update the selected cell as recorded in the sheet view
$ws.Range("G25").Select()
